$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header cell H1 into the two new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the new header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set the new data values (plain, matching H2's unstyled format)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
